$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to match the Payments format
$ws.Range("D1").Value = "ACC."
$ws.Range("E1").Value = "REC."
$ws.Range("J1").Value = "FRT."
